$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the URL in A2 and mark "Privacy Found" as True
$ws.Range("A2").Value = "http://agrar.bayer.de/"
$ws.Range("B2").Value = $true

# Remove the duplicate data row (row 3)
$ws.Rows("3").Delete()

# Restore the normal font size (11) for the header row style
$ws.Range("A1:H1").Font.Size = 11

# Auto-fit the columns whose content/font actually changed so their
# widths reflect the new text (column C - "Privacy Name" - is untouched
# since neither its header nor its data changed)
foreach ($colLetter in @("A", "B", "D", "E", "F", "G", "H")) {
    $col = $ws.Columns($colLetter)
    $col.AutoFit()
    $w = $col.ColumnWidth
    $col.ColumnWidth = $w
}
